$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.151.98"
$ws.Range("E2").Value = "  -0.92%  "

$ws.Range("D3").Value = "3.318.86"
$ws.Range("E3").Value = "  -1.10%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.54"
$ws.Range("E5").Value = "  +2.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.85"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.653"
$ws.Range("E7").Value = "  +4.17%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "3.318.21"
$ws.Range("E9").Value = "  -1.08%  "

$ws.Range("E10").Value = "  -2.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.81"
$ws.Range("E11").Value = "  +2.41%  "

$ws.Range("E12").Value = "  -0.36%  "

$ws.Range("D13").Value = "3.905.50"
$ws.Range("E13").Value = "  -0.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.132"
$ws.Range("E14").Value = "  -2.44%  "

$ws.Range("D15").Value = "66.208.30"
$ws.Range("E15").Value = "  -0.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.31"
$ws.Range("E16").Value = "  -2.26%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000165"
$ws.Range("E17").Value = "  -1.43%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.314.25"
$ws.Range("E18").Value = "  -1.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "424.33"
$ws.Range("E19").Value = "  -3.80%  "

$ws.Range("E20").Value = "  -2.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.17"
$ws.Range("E21").Value = "  -3.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.39"
$ws.Range("E22").Value = "  -2.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.85"
$ws.Range("E23").Value = "  -2.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("D26").Value = "3.464.14"
$ws.Range("E26").Value = "  -0.89%  "

$ws.Range("E27").Value = "  -1.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.203"
$ws.Range("E28").Value = "  +4.71%  "

$ws.Range("E29").Value = "  -2.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.98"
$ws.Range("E30").Value = "  -1.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.92"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.40"
$ws.Range("E33").Value = "  -2.01%  "

$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("E35").Value = "  -2.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.59"
$ws.Range("E36").Value = "  -3.17%  "

$ws.Range("E37").Value = "  -3.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.67"
$ws.Range("E38").Value = "  -0.58%  "

$ws.Range("E39").Value = "  -2.96%  "

$ws.Range("D40").Value = "2.884.11"
$ws.Range("E40").Value = "  +1.53%  "

$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.38"
$ws.Range("E42").Value = "  -5.32%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.32"
$ws.Range("E43").Value = "  -2.68%  "

$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.759"
$ws.Range("E44").Value = "  -5.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.80"
$ws.Range("E45").Value = "  -1.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0662"
$ws.Range("E46").Value = "  -0.95%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.92"
$ws.Range("E47").Value = "  -4.57%  "

$ws.Range("E48").Value = "  -1.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.11"
$ws.Range("E49").Value = "  -5.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "313.54"
$ws.Range("E50").Value = "  -3.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0273"
$ws.Range("E51").Value = "  -0.05%  "
